$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 419, shifting old rows 419:500 down to 420:501
$ws.Rows.Item(419).Insert()

# Populate the newly inserted row 419 with a new Zanahoria record
$ws.Range("A419").Value = 5
$ws.Range("B419").Value = "Macroferia Regional de Talca"
$ws.Range("C419").Value = "Maule"
$ws.Range("D419").Value = 45015
$ws.Range("E419").Value = 7
$ws.Range("F419").Value = 100114013
$ws.Range("G419").Value = "Zanahoria"
$ws.Range("H419").Value = "Sin especificar"
$ws.Range("I419").Value = "Primera"
$ws.Range("J419").Value = 600
$ws.Range("K419").Value = 7000
$ws.Range("L419").Value = 7000
$ws.Range("M419").Value = 7000
$ws.Range("N419").Value = "`$/saco 20 kilos"
$ws.Range("O419").Value = "Región de Ñuble"
$ws.Range("P419").Value = 350
$ws.Range("Q419").Value = 20
$ws.Range("R419").Value = "Hortaliza"
